{"js": "// Office.js (Word JavaScript API) script\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\n\n// 1) Remove the old \"_GoBack\" bookmark that currently sits alone in the\n//    trailing empty paragraph at the end of the document (it will be\n//    re-created at the top of the title paragraph below). Bookmark names\n//    must be unique, so delete the old one first.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Load the paragraphs so we can grab the title paragraph (first one).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\nconst startRange = titlePara.getRange(\"Start\");\n\n// 3) Insert \" Q\", a tab run, and a new \"_GoBack\" bookmark at the very\n//    start of the title paragraph, before the existing \"Billeder til \"\n//    run. We use insertOoxml so we get a real <w:tab/> run (not a literal\n//    tab character) and bookmark elements, matching native Word output.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\"> Q</w:t></w:r>' +\n  '<w:r><w:tab/></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nstartRange.insertOoxml(ooxml, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# $word / $d (ActiveDocument) are pre-seeded by the host.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the old \"_GoBack\" bookmark that currently sits alone in the\n#    trailing empty paragraph at the end of the document (it will be\n#    re-created at the top of the title paragraph below). Bookmark names\n#    must stay unique, so delete the old one first. \"_GoBack\" is a hidden\n#    bookmark, so address it by name rather than relying on it showing up\n#    while enumerating Bookmarks.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Insert \" Q\", a tab run, and a new \"_GoBack\" bookmark at the very\n#    start of the title paragraph (the first paragraph), before the\n#    existing \"Billeder til \" run. InsertXML gives us a real <w:tab/>\n#    run (not a literal tab character) and bookmark elements, matching\n#    native Word output.\n$p1 = $d.Paragraphs(1)\n$r = $p1.Range\n$r.Collapse(1)  # wdCollapseStart\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n  '<pkg:xmlData>' + `\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n  '<w:body><w:p>' + `\n  '<w:r><w:t xml:space=\"preserve\"> Q</w:t></w:r>' + `\n  '<w:r><w:tab/></w:r>' + `\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' + `\n  '</w:p></w:body></w:document>' + `\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$r.InsertXML($xml)\n"}
